# Update gh-pages to output generated at 456a3b4
# Apply updated "want to go" / ticket counts to the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value = 45
$ws1.Range("F5").Value = 1709
$ws1.Range("F7").Value = 2164
$ws1.Range("F11").Value = 4835
$ws1.Range("F16").Value = 29
$ws1.Range("F21").Value = 3763
$ws1.Range("F23").Value = 630
$ws1.Range("F27").Value = 115
$ws1.Range("F34").Value = 884
$ws1.Range("F35").Value = 2398

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value = 45
$ws4.Range("F5").Value = 1709
$ws4.Range("F7").Value = 2164
$ws4.Range("F11").Value = 4835
$ws4.Range("F16").Value = 29
$ws4.Range("F21").Value = 3763
$ws4.Range("F23").Value = 630
$ws4.Range("F27").Value = 115
$ws4.Range("F35").Value = 884
$ws4.Range("F36").Value = 2398
